# Generate Report for Handback
# Adds a new "row 4" entry (file b086b3ff-5c41-4304-82cf-7f905841a259.md) to the
# Overview / zh-cn / de-de sheets, mirroring the existing row layout, extends the
# backing tables, and wires up the corresponding hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"  (A1:G3 -> A1:G4)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A4").Value = "b086b3ff-5c41-4304-82cf-7f905841a259.md"
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee8d97313c08797298157e08e877d1f7fe878e38/e2e/b086b3ff-5c41-4304-82cf-7f905841a259.md", "", "", "e2e\b086b3ff-5c41-4304-82cf-7f905841a259.md")
$ws1.Range("B4").Style = "HyperLink"
$ws1.Range("C4").Value = ".md"
$ws1.Range("E4").Value = "Handed back: in sync with en-US"
$ws1.Range("F4").Value = "Handed back: in sync with en-US"
$ws1.Range("G4").Value = "2016-09-04 22:48:51"
$ws1.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"  (A1:P3 -> A1:P4)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7cac617eaee3e01ab991ad5f8200e1dfe0309f64/e2e/b086b3ff-5c41-4304-82cf-7f905841a259.md", "", "", "b086b3ff-5c41-4304-82cf-7f905841a259.md")
$ws2.Range("A4").Style = "HyperLink"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Handed back: in sync with en-US"
$ws2.Range("D4").Value = "e2e"
$ws2.Range("E4").Value = "ht"
$ws2.Range("F4").Value = "'True"
$ws2.Range("G4").Value = "b086b3ff-5c41-4304-82cf-7f905841a259.964570261333555f8fc5f0a155950fdb5a74951a.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-09-04 22:48:46"
$ws2.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Hyperlinks.Add($ws2.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7cac617eaee3e01ab991ad5f8200e1dfe0309f64/e2e/b086b3ff-5c41-4304-82cf-7f905841a259.md", "", "", "b086b3ff-5c41-4304-82cf-7f905841a259.md")
$ws2.Range("I4").Style = "HyperLink"
$ws2.Range("J4").Value = "b086b3ff-5c41-4304-82cf-7f905841a259.964570261333555f8fc5f0a155950fdb5a74951a.zh-cn.xlf"
$ws2.Range("K4").Value = "2016-09-04 22:49:07"
$ws2.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("M4").Value = "'True"
$ws2.Range("O4").Value = "'False"

$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de"  (A1:P3 -> A1:P4)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/666608268beb068a4e1dd2a4d3c058a6719e24f6/e2e/b086b3ff-5c41-4304-82cf-7f905841a259.md", "", "", "b086b3ff-5c41-4304-82cf-7f905841a259.md")
$ws3.Range("A4").Style = "HyperLink"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Handed back: in sync with en-US"
$ws3.Range("D4").Value = "e2e"
$ws3.Range("E4").Value = "ht"
$ws3.Range("F4").Value = "'True"
$ws3.Range("G4").Value = "b086b3ff-5c41-4304-82cf-7f905841a259.964570261333555f8fc5f0a155950fdb5a74951a.de-de.xlf"
$ws3.Range("H4").Value = "2016-09-04 22:48:51"
$ws3.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Hyperlinks.Add($ws3.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/666608268beb068a4e1dd2a4d3c058a6719e24f6/e2e/b086b3ff-5c41-4304-82cf-7f905841a259.md", "", "", "b086b3ff-5c41-4304-82cf-7f905841a259.md")
$ws3.Range("I4").Style = "HyperLink"
$ws3.Range("J4").Value = "b086b3ff-5c41-4304-82cf-7f905841a259.964570261333555f8fc5f0a155950fdb5a74951a.de-de.xlf"
$ws3.Range("K4").Value = "2016-09-04 22:49:15"
$ws3.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("M4").Value = "'True"
$ws3.Range("O4").Value = "'False"

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:P4"))

Write-Output "Generate Report for Handback: done"
